# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt -
# Cilantro" above the current row 535. This pushes the existing rows
# 535..547 down to 536..548 (dimension grows from A1:R547 to A1:R548) and
# fills the freshly inserted row 535 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 535:547 down by one, duplicating row 535's formatting into the
# newly created (blank) row 535 - matches Excel's default "Insert" behavior.
$ws.Rows.Item(535).Insert()

$ws.Range("A535").Value = 4
$ws.Range("B535").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C535").Value = "Los Lagos"
$ws.Range("D535").Value = 45239
$ws.Range("E535").Value = 10
$ws.Range("F535").Value = 100112040
$ws.Range("G535").Value = "Cilantro"
$ws.Range("H535").Value = "Sin especificar"
$ws.Range("I535").Value = "Primera"
$ws.Range("J535").Value = 80
$ws.Range("K535").Value = 13000
$ws.Range("L535").Value = 13000
$ws.Range("M535").Value = 13000
$ws.Range("N535").Value = "$/caja 36 atados"
$ws.Range("O535").Value = "Región Metropolitana"
$ws.Range("P535").Value = 361
$ws.Range("Q535").Value = 36
$ws.Range("R535").Value = "Hortaliza"
